# Work order templates cleanup:
#  - remove the stray "8/16/2021" test date from the WO receipt traveler sheet
#  - refresh the cursor/selection on each sheet, leaving "WO receipt" as the
#    active (front) sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Create WO
$ws2 = $wb.Worksheets.Item(2)   # Time and Qty Booking
$ws3 = $wb.Worksheets.Item(3)   # WO receipt

# Clear the leftover test value in the WO receipt traveler sheet
$ws3.Range("C2").ClearContents() | Out-Null

# Move the selection on "Create WO" to F10
$ws1.Range("F10").Select() | Out-Null

# Finish on "WO receipt", making it the active sheet, with C2 selected
$ws3.Activate() | Out-Null
$ws3.Range("C2").Select() | Out-Null
